# Update view-count-like figures (column F) on the "展览" and "全部类型"
# sheets to reflect newly generated data (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12874
$ws1.Range("F10").Value = 12794
$ws1.Range("F13").Value = 8662
$ws1.Range("F14").Value = 7664
$ws1.Range("F19").Value = 979

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12874
$ws4.Range("F11").Value = 12794
$ws4.Range("F14").Value = 8662
$ws4.Range("F15").Value = 7664
$ws4.Range("F20").Value = 979
